# Slide 9 already has a semi-transparent red (FF0000) highlight rectangle
# that only covers bit-index columns 26-31 of the top bit table, leaving a
# gap at column 25 before the adjacent yellow highlight begins. Fix the
# color-coding by adding a new red highlight rectangle, built from the
# existing one (so it inherits the exact same fill/shape style), sized to
# also cover the missing column.
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)

# The existing red (FF0000) highlight rectangle on this slide.
$src = $s.Shapes.Item(10)

$dup = $src.Duplicate()
$newShape = $dup.Item(1)
$newShape.Name = ""

# Exact target geometry (EMU): off x=769358 y=2957677, ext cx=2418729 cy=316926
$newShape.Left = 60.57938007874016
$newShape.Top = 232.8879527559055
$newShape.Width = 190.45111236220473
$newShape.Height = 24.9548031496063
